# Actualización automática del inventario: agrega un nuevo producto (fila 31)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31

$ws.Cells.Item($row, 1).Value = "LJNL6J"
$ws.Cells.Item($row, 2).Value = "Drum DL-410 para Tambor de Unidad de imagen Pantum"
$ws.Cells.Item($row, 3).Value = "P3010DW P3300DN P3300DW M6700DW M6800FDW M7100DN M7200FDW M7300FDW M7300FDN"
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 100000
$ws.Cells.Item($row, 6).Value = 9
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Formula = "=(E31-D31)*G31"
$ws.Cells.Item($row, 9).Formula = "=D31*F31"
$ws.Cells.Item($row, 10).Value = 0
